# Auto-generated Excel COM-interop script to update paises.xlsx
# Applies: (1) updated "Datos actualizados" timestamp, (2) reordered country
# rows (new entrants: Egipto, Costa Rica, Costa de Marfil) causing rank shifts,
# and (3) refreshed case-count statistics for many countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 20:16"

# Full country table (row 4 - row 199): Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    ,@("China",81171,78,73159,4735,1573,7,3277)
    ,@("Italia",69176,5249,8326,54030,3393,743,6820)
    ,@("Estados Unidos",51678,7944,368,50643,1175,114,667)
    ,@("España",39676,4540,3794,33082,2355,489,2800)
    ,@("Alemania",32781,3725,3243,29382,23,33,156)
    ,@("Iran",24811,1762,8913,13964,0,122,1934)
    ,@("Francia",22304,2448,3281,17923,2516,240,1100)
    ,@("Suiza",9877,1082,131,9624,141,2,122)
    ,@("Corea del Sur",9037,76,3507,5410,59,9,120)
    ,@("Reino Unido",8077,1427,135,7520,20,87,422)
    ,@("Paises Bajos",5560,811,2,5282,435,63,276)
    ,@("Austria",5137,663,9,5100,22,7,28)
    ,@("Belgica",4269,526,461,3686,381,34,122)
    ,@("Noruega",2779,154,6,2761,44,2,12)
    ,@("Canada",2583,492,112,2446,1,1,25)
    ,@("Portugal",2362,302,22,2307,48,10,33)
    ,@("Suecia",2286,240,16,2234,136,9,36)
    ,@("Australia",2144,257,118,2018,11,1,8)
    ,@("Brasil",1980,56,2,1944,18,0,34)
    ,@("Israel",1656,214,49,1604,31,2,3)
    ,@("Malasia",1624,106,183,1425,64,2,16)
    ,@("Dinamarca",1591,131,1,1558,69,8,32)
    ,@("Turquia",1529,0,0,1492,0,0,37)
    ,@("Chequia",1394,158,8,1384,19,1,2)
    ,@("Japon",1193,65,285,865,54,1,43)
    ,@("Irlanda",1125,0,5,1114,29,0,6)
    ,@("Luxemburgo",1099,224,6,1085,3,0,8)
    ,@("Ecuador",1049,68,3,1019,2,9,27)
    ,@("Pakistan",972,97,18,947,5,1,7)
    ,@("Chile",922,176,17,903,7,0,2)
    ,@("Polonia",884,135,1,873,3,2,10)
    ,@("Tailandia",827,106,52,771,7,3,4)
    ,@("Finlandia",792,92,10,781,11,0,1)
    ,@("Arabia Saudita",767,205,28,738,0,1,1)
    ,@("Rumania",762,186,79,675,15,1,8)
    ,@("Grecia",743,48,29,694,35,3,20)
    ,@("Crucero",712,0,587,115,15,2,10)
    ,@("Indonesia",686,107,30,601,0,6,55)
    ,@("Islandia",648,60,51,595,13,1,2)
    ,@("Singapur",558,49,156,400,14,0,2)
    ,@("Sudafrica",554,152,4,550,2,0,0)
    ,@("Filipinas",552,90,20,497,1,2,35)
    ,@("India",536,37,40,486,0,0,10)
    ,@("Catar",526,25,41,485,6,0,0)
    ,@("Rusia",495,57,22,472,8,0,1)
    ,@("Eslovenia",480,38,3,473,12,1,4)
    ,@("Peru",416,21,1,408,19,2,7)
    ,@("Egipto",402,36,80,302,0,1,20)
    ,@("Barein",392,15,177,212,2,1,3)
    ,@("Hong Kong",386,29,102,280,4,0,4)
    ,@("Croacia",382,67,5,376,6,0,1)
    ,@("Estonia",369,17,7,362,4,0,0)
    ,@("Mexico",367,51,4,359,1,1,4)
    ,@("Panama",345,0,1,338,33,0,6)
    ,@("Libano",318,51,8,306,4,0,4)
    ,@("Irak",316,50,75,214,0,4,27)
    ,@("Republica Dominicana",312,67,3,303,0,3,6)
    ,@("Colombia",306,29,6,297,0,0,3)
    ,@("Serbia",303,54,15,285,21,0,3)
    ,@("Argentina",301,0,52,243,0,2,6)
    ,@("Argelia",264,34,24,221,0,2,19)
    ,@("Armenia",249,14,4,245,6,0,0)
    ,@("Emiratos Arabes Unidos",248,50,45,201,2,0,2)
    ,@("Bulgaria",218,17,3,212,8,0,3)
    ,@("Taiwan",216,21,29,185,0,0,2)
    ,@("Lituania",209,30,1,207,1,0,1)
    ,@("Eslovaquia",204,18,7,197,2,0,0)
    ,@("Letonia",197,17,1,196,0,0,0)
    ,@("Kuwait",191,2,39,152,5,0,0)
    ,@("San Marino",187,0,4,162,12,1,21)
    ,@("Hungria",187,20,21,157,6,1,9)
    ,@("Costa Rica",177,19,2,173,4,0,2)
    ,@("Marruecos",170,27,6,159,1,1,5)
    ,@("Principado de Andorra",164,31,1,162,7,0,1)
    ,@("Bosnia y Herzegovina",164,28,2,160,1,1,2)
    ,@("Uruguay",162,0,0,162,3,0,0)
    ,@("Nueva Zelanda",155,53,12,143,0,0,0)
    ,@("Jordania",153,26,1,152,0,0,0)
    ,@("Republica de Macedonia",148,12,1,145,1,0,2)
    ,@("Vietnam",134,11,17,117,3,0,0)
    ,@("Moldavia",125,16,2,122,10,0,1)
    ,@("Republica de Chipre",124,8,3,118,3,2,3)
    ,@("Albania",123,19,10,108,2,1,5)
    ,@("Islas Feroe",122,4,23,99,0,0,0)
    ,@("Tunez",114,25,1,110,11,0,3)
    ,@("Burkina Faso",114,15,7,103,0,0,4)
    ,@("Malta",110,3,2,108,1,0,0)
    ,@("Brunei",104,13,2,102,2,0,0)
    ,@("Sri Lanka",102,5,2,100,2,0,0)
    ,@("Ucrania",97,24,1,93,0,0,3)
    ,@("Camboya",91,4,4,87,1,0,0)
    ,@("Azerbaiyan",87,15,10,76,6,0,1)
    ,@("Senegal",86,7,8,78,0,0,0)
    ,@("Venezuela",84,0,15,69,2,0,0)
    ,@("Oman",84,18,17,67,0,0,0)
    ,@("Reunion",83,12,1,82,0,0,0)
    ,@("Bielorrusia",81,0,22,59,0,0,0)
    ,@("Afganistan",74,34,1,72,0,0,1)
    ,@("Costa de Marfil",73,48,2,71,0,0,0)
    ,@("Kazajistan",72,10,0,72,0,0,0)
    ,@("Georgia",70,9,9,61,1,0,0)
    ,@("Camerun",66,10,2,64,0,0,0)
    ,@("Guadalupe",62,0,0,61,4,0,1)
    ,@("Estado de Palestina",60,1,16,44,0,0,0)
    ,@("Trinidad yTobago",53,2,0,53,0,0,0)
    ,@("Martinica",53,0,0,52,7,0,1)
    ,@("Ghana",52,25,0,50,0,0,2)
    ,@("Liechtenstein",51,0,0,51,0,0,0)
    ,@("Uzbekistan",50,4,0,50,0,0,0)
    ,@("Cuba",48,8,0,47,2,0,1)
    ,@("Montenegro",47,20,0,46,0,0,1)
    ,@("Consejo Danes para los Refugiados",45,9,0,43,0,0,2)
    ,@("Nigeria",44,4,2,41,0,0,1)
    ,@("Kirguistan",42,26,0,42,0,0,0)
    ,@("Mauricio",42,6,0,40,1,0,2)
    ,@("Puerto Rico",39,8,1,36,0,0,2)
    ,@("Banglades",39,6,5,30,0,1,4)
    ,@("Mayotte",36,12,0,36,0,0,0)
    ,@("Ruanda",36,0,0,36,0,0,0)
    ,@("Guam",32,3,0,31,0,0,1)
    ,@("Honduras",30,0,0,30,0,0,0)
    ,@("Bolivia",29,2,0,29,0,0,0)
    ,@("Macao",28,3,10,18,0,0,0)
    ,@("Paraguay",27,5,0,25,1,1,2)
    ,@("Kenia",25,9,0,25,0,0,0)
    ,@("Polinesia Francesa",23,5,0,23,0,0,0)
    ,@("Isla de Man",23,10,0,23,0,0,0)
    ,@("Monaco",23,0,1,22,0,0,0)
    ,@("Guayana Francesa",23,3,6,17,0,0,0)
    ,@("Guatemala",21,1,0,20,0,0,1)
    ,@("Jamaica",21,2,2,18,0,0,1)
    ,@("Togo",20,2,1,19,0,0,0)
    ,@("Guyana",20,0,0,19,0,0,1)
    ,@("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0)
    ,@("Madagascar",17,5,0,17,0,0,0)
    ,@("Barbados",17,0,0,17,0,0,0)
    ,@("Gibraltar",15,0,5,10,0,0,0)
    ,@("Maldivas",13,0,5,8,0,0,0)
    ,@("Etiopia",12,1,0,12,0,0,0)
    ,@("Tanzania",12,0,0,12,0,0,0)
    ,@("Aruba",12,3,1,11,0,0,0)
    ,@("Mongolia",10,0,0,10,0,0,0)
    ,@("Nueva Caledonia",10,2,0,10,0,0,0)
    ,@("Guinea Ecuatorial",9,0,0,9,0,0,0)
    ,@("Uganda",9,0,0,9,0,0,0)
    ,@("San Martin (Parte Francesa)",8,0,0,8,0,0,0)
    ,@("Seychelles",7,0,0,7,0,0,0)
    ,@("Bermudas",6,0,0,6,0,0,0)
    ,@("Haiti",6,0,0,6,0,0,0)
    ,@("Surinam",6,1,0,6,0,0,0)
    ,@("Benin",6,0,0,6,0,0,0)
    ,@("Gabon",6,0,0,5,0,0,1)
    ,@("Islas Caimanes",6,1,0,5,0,0,1)
    ,@("Namibia",6,2,2,4,0,0,0)
    ,@("El Salvador",5,2,0,5,0,0,0)
    ,@("Groenlandia",5,1,2,3,0,0,0)
    ,@("Guinea",4,0,0,4,0,0,0)
    ,@("Congo",4,0,0,4,0,0,0)
    ,@("Bahamas",4,0,0,4,0,0,0)
    ,@("Santa Sede",4,3,0,4,0,0,0)
    ,@("Suazilandia",4,0,0,4,0,0,0)
    ,@("Fiyi",4,1,0,4,0,0,0)
    ,@("Curazao",4,0,0,3,0,0,1)
    ,@("Republica del Chad",3,1,0,3,0,0,0)
    ,@("Mozambique",3,2,0,3,0,0,0)
    ,@("Antigua y Barbuda",3,0,0,3,0,0,0)
    ,@("Liberia",3,0,0,3,0,0,0)
    ,@("Angola",3,0,0,3,0,0,0)
    ,@("Niger",3,0,0,3,0,0,0)
    ,@("Santa Lucia",3,0,0,3,0,0,0)
    ,@("Zambia",3,0,0,3,0,0,0)
    ,@("Republica de Africa Central",3,0,0,3,0,0,0)
    ,@("Republica de Yibuti",3,0,0,3,0,0,0)
    ,@("San Bartolome",3,0,0,3,0,0,0)
    ,@("Sudan",3,1,0,2,0,0,1)
    ,@("Gambia",3,1,0,2,0,0,1)
    ,@("Cabo Verde",3,0,0,2,0,1,1)
    ,@("Zimbabue",3,0,0,2,0,0,1)
    ,@("Mauritania",2,0,0,2,0,0,0)
    ,@("Butan",2,0,0,2,0,0,0)
    ,@("Laos",2,2,0,2,0,0,0)
    ,@("San Martin (Parte Holandesa)",2,0,0,2,0,0,0)
    ,@("Dominica",2,0,0,2,0,0,0)
    ,@("Birmania",2,0,0,2,0,0,0)
    ,@("Nicaragua",2,0,0,2,0,0,0)
    ,@("Nepal",2,0,1,1,0,0,0)
    ,@("Belice",1,0,0,1,0,0,0)
    ,@("Papua Nueva Guinea",1,0,0,1,0,0,0)
    ,@("Timor Oriental",1,0,0,1,0,0,0)
    ,@("Eritrea",1,0,0,1,0,0,0)
    ,@("Siria",1,0,0,1,0,0,0)
    ,@("Granada",1,0,0,1,0,0,0)
    ,@("Islas Turcas y Caicos",1,0,0,1,0,0,0)
    ,@("Montserrat",1,0,0,1,0,0,0)
    ,@("Somalia",1,0,0,1,0,0,0)
    ,@("San Vicente y las Granadinas",1,0,0,1,0,0,0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 4 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = [double]$rec[1]
    $ws.Cells.Item($row, 3).Value = [double]$rec[2]
    $ws.Cells.Item($row, 4).Value = [double]$rec[3]
    $ws.Cells.Item($row, 5).Value = [double]$rec[4]
    $ws.Cells.Item($row, 6).Value = [double]$rec[5]
    $ws.Cells.Item($row, 7).Value = [double]$rec[6]
    $ws.Cells.Item($row, 8).Value = [double]$rec[7]
}

Write-Host "Applied $($data.Count) row updates."
